$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# NOTE: "11×75=" is both a target value (from 83×39=) and a source value
# (original cell that becomes 43×51=). The original "11×75=" occurs later
# in the document than "83×39=", so we must replace it BEFORE turning
# "83×39=" into "11×75=", otherwise the newly created "11×75=" text would
# also get matched and incorrectly turned into "43×51=".
Replace-Text "11×75=" "43×51="

Replace-Text "37×42=" "23×51="
Replace-Text "85×14=" "24×13="
Replace-Text "96×48=" "85×74="
Replace-Text "64×40=" "75×63="
Replace-Text "77×43=" "69×66="
Replace-Text "70×82=" "18×82="
Replace-Text "40×59=" "59×95="
Replace-Text "67×45=" "21×11="
Replace-Text "45×51=" "32×90="
Replace-Text "78×17=" "90×61="
Replace-Text "42×46=" "74×49="
Replace-Text "83×39=" "11×75="
Replace-Text "48×52=" "62×70="
Replace-Text "86×96=" "29×95="
Replace-Text "39×12=" "20×42="
Replace-Text "49×21=" "59×13="
Replace-Text "18×31=" "37×81="
Replace-Text "98×94=" "88×17="
Replace-Text "21×16=" "62×76="
Replace-Text "27×94=" "70×56="
Replace-Text "87×52=" "58×58="
Replace-Text "20×39=" "57×19="
Replace-Text "96×59=" "60×33="
Replace-Text "12×22=" "78×36="
